# Update the "Förändrad" (changed) date in column C for the data rows
# (rows 2-41) from 2023-09-14 (serial 45183) to 2023-09-15 (serial 45184).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 41; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
